# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> was "Office Theme" (used only by the Notes Master)
#   ppt/theme/theme2.xml  -> was "Integral"     (used by the Slide Master / all slides)
#
# The authored edit swaps the two themes' content: the slide-facing theme
# (theme2.xml) becomes the stock "Office Theme" color palette, while the
# notes-only theme (theme1.xml) becomes the "Integral" palette.
#
# The PowerPoint object model only exposes a writable per-channel RGB on the
# color scheme that backs the active Slide Master (Slide.ThemeColorScheme /
# SlideRange.ThemeColorScheme), which is the visually-significant half of
# the swap. Drive that color-by-color to the Office Theme values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> (theme slot, target "Office Theme" RGB)
# RGB() packing used by PowerPoint's COM automation is R + G*256 + B*65536.
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
